$wb = $excel.ActiveWorkbook

# ===== Station1: update carry-demand values =====
$ws = $wb.Worksheets.Item("Station1")
$ws.Range("D2").Value = 182
$ws.Range("E2").Value = 183
$ws.Range("F2").Value = 184
$ws.Range("G2").Value = 185
$ws.Range("D3").Value = 192
$ws.Range("E3").Value = 193
$ws.Range("F3").Value = 194
$ws.Range("G3").Value = 195
$ws.Range("D4").Value = 102
$ws.Range("E4").Value = 103
$ws.Range("F4").Value = 104
$ws.Range("G4").Value = 105
$ws.Range("D5").Value = 112
$ws.Range("E5").Value = 113
$ws.Range("F5").Value = 114
$ws.Range("G5").Value = 115
$ws.Range("D6").Value = 122
$ws.Range("E6").Value = 123
$ws.Range("F6").Value = 124
$ws.Range("G6").Value = 125
$ws.Range("D7").Value = 132
$ws.Range("E7").Value = 133
$ws.Range("F7").Value = 134
$ws.Range("G7").Value = 135
$ws.Range("D8").Value = 142
$ws.Range("E8").Value = 143
$ws.Range("F8").Value = 144
$ws.Range("G8").Value = 145
$ws.Range("D9").Value = 152
$ws.Range("E9").Value = 153
$ws.Range("F9").Value = 154
$ws.Range("G9").Value = 155
$ws.Range("D10").Value = 162
$ws.Range("E10").Value = 163
$ws.Range("F10").Value = 164
$ws.Range("G10").Value = 165
$ws.Range("D11").Value = 172
$ws.Range("E11").Value = 173
$ws.Range("F11").Value = 174
$ws.Range("G11").Value = 175

# ===== Station2: update carry-demand values =====
$ws = $wb.Worksheets.Item("Station2")
$ws.Range("C2").Value = 281
$ws.Range("E2").Value = 283
$ws.Range("F2").Value = 284
$ws.Range("G2").Value = 285
$ws.Range("C3").Value = 291
$ws.Range("E3").Value = 293
$ws.Range("F3").Value = 294
$ws.Range("G3").Value = 295
$ws.Range("C4").Value = 201
$ws.Range("E4").Value = 203
$ws.Range("F4").Value = 204
$ws.Range("G4").Value = 205
$ws.Range("C5").Value = 211
$ws.Range("E5").Value = 213
$ws.Range("F5").Value = 214
$ws.Range("G5").Value = 215
$ws.Range("C6").Value = 221
$ws.Range("E6").Value = 223
$ws.Range("F6").Value = 224
$ws.Range("G6").Value = 225
$ws.Range("C7").Value = 231
$ws.Range("E7").Value = 233
$ws.Range("F7").Value = 234
$ws.Range("G7").Value = 235
$ws.Range("C8").Value = 241
$ws.Range("E8").Value = 243
$ws.Range("F8").Value = 244
$ws.Range("G8").Value = 245
$ws.Range("C9").Value = 251
$ws.Range("E9").Value = 253
$ws.Range("F9").Value = 254
$ws.Range("G9").Value = 255
$ws.Range("C10").Value = 261
$ws.Range("E10").Value = 263
$ws.Range("F10").Value = 264
$ws.Range("G10").Value = 265
$ws.Range("C11").Value = 271
$ws.Range("E11").Value = 273
$ws.Range("F11").Value = 274
$ws.Range("G11").Value = 275

# ===== Station3: update carry-demand values =====
$ws = $wb.Worksheets.Item("Station3")
$ws.Range("C2").Value = 381
$ws.Range("D2").Value = 382
$ws.Range("F2").Value = 384
$ws.Range("G2").Value = 385
$ws.Range("C3").Value = 391
$ws.Range("D3").Value = 392
$ws.Range("F3").Value = 394
$ws.Range("G3").Value = 395
$ws.Range("C4").Value = 301
$ws.Range("D4").Value = 302
$ws.Range("F4").Value = 304
$ws.Range("G4").Value = 305
$ws.Range("C5").Value = 311
$ws.Range("D5").Value = 312
$ws.Range("F5").Value = 314
$ws.Range("G5").Value = 315
$ws.Range("C6").Value = 321
$ws.Range("D6").Value = 322
$ws.Range("F6").Value = 324
$ws.Range("G6").Value = 325
$ws.Range("C7").Value = 331
$ws.Range("D7").Value = 332
$ws.Range("F7").Value = 334
$ws.Range("G7").Value = 335
$ws.Range("C8").Value = 341
$ws.Range("D8").Value = 342
$ws.Range("F8").Value = 344
$ws.Range("G8").Value = 345
$ws.Range("C9").Value = 351
$ws.Range("D9").Value = 352
$ws.Range("F9").Value = 354
$ws.Range("G9").Value = 355
$ws.Range("C10").Value = 361
$ws.Range("D10").Value = 362
$ws.Range("F10").Value = 364
$ws.Range("G10").Value = 365
$ws.Range("C11").Value = 371
$ws.Range("D11").Value = 372
$ws.Range("F11").Value = 374
$ws.Range("G11").Value = 375
# restore/align formatting for Station3
$ws.Range("F2").Font.Color = 0
$ws.Range("G2").Font.Color = 0
$ws.Range("F3").Font.Color = 0
$ws.Range("G3").Font.Color = 0
$ws.Range("F4").Font.Color = 0
$ws.Range("G4").Font.Color = 0
$ws.Range("F5").Font.Color = 0
$ws.Range("G5").Font.Color = 0
$ws.Range("F6").Font.Color = 0
$ws.Range("G6").Font.Color = 0
$ws.Range("F7").Font.Color = 0
$ws.Range("G7").Font.Color = 0
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Style = "Normal"

# ===== Station4: update carry-demand values =====
$ws = $wb.Worksheets.Item("Station4")
$ws.Range("C2").Value = 481
$ws.Range("D2").Value = 482
$ws.Range("E2").Value = 483
$ws.Range("G2").Value = 485
$ws.Range("C3").Value = 491
$ws.Range("D3").Value = 492
$ws.Range("E3").Value = 493
$ws.Range("G3").Value = 495
$ws.Range("C4").Value = 401
$ws.Range("D4").Value = 402
$ws.Range("E4").Value = 403
$ws.Range("G4").Value = 405
$ws.Range("C5").Value = 411
$ws.Range("D5").Value = 412
$ws.Range("E5").Value = 413
$ws.Range("G5").Value = 415
$ws.Range("C6").Value = 421
$ws.Range("D6").Value = 422
$ws.Range("E6").Value = 423
$ws.Range("G6").Value = 425
$ws.Range("C7").Value = 431
$ws.Range("D7").Value = 432
$ws.Range("E7").Value = 433
$ws.Range("G7").Value = 435
$ws.Range("C8").Value = 441
$ws.Range("D8").Value = 442
$ws.Range("E8").Value = 443
$ws.Range("G8").Value = 445
$ws.Range("C9").Value = 451
$ws.Range("D9").Value = 452
$ws.Range("E9").Value = 453
$ws.Range("G9").Value = 455
$ws.Range("C10").Value = 461
$ws.Range("D10").Value = 462
$ws.Range("E10").Value = 463
$ws.Range("G10").Value = 465
$ws.Range("C11").Value = 471
$ws.Range("D11").Value = 472
$ws.Range("E11").Value = 473
$ws.Range("G11").Value = 475
# restore/align formatting for Station4
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Style = "Normal"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Style = "Normal"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Style = "Normal"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Style = "Normal"

# ===== Station5: update carry-demand values =====
$ws = $wb.Worksheets.Item("Station5")
$ws.Range("C2").Value = 581
$ws.Range("D2").Value = 582
$ws.Range("E2").Value = 583
$ws.Range("F2").Value = 584
$ws.Range("C3").Value = 591
$ws.Range("D3").Value = 592
$ws.Range("E3").Value = 593
$ws.Range("F3").Value = 594
$ws.Range("C4").Value = 501
$ws.Range("D4").Value = 502
$ws.Range("E4").Value = 503
$ws.Range("F4").Value = 504
$ws.Range("C5").Value = 511
$ws.Range("D5").Value = 512
$ws.Range("E5").Value = 513
$ws.Range("F5").Value = 514
$ws.Range("C6").Value = 521
$ws.Range("D6").Value = 522
$ws.Range("E6").Value = 523
$ws.Range("F6").Value = 524
$ws.Range("C7").Value = 531
$ws.Range("D7").Value = 532
$ws.Range("E7").Value = 533
$ws.Range("F7").Value = 534
$ws.Range("C8").Value = 541
$ws.Range("D8").Value = 542
$ws.Range("E8").Value = 543
$ws.Range("F8").Value = 544
$ws.Range("C9").Value = 551
$ws.Range("D9").Value = 552
$ws.Range("E9").Value = 553
$ws.Range("F9").Value = 554
$ws.Range("C10").Value = 561
$ws.Range("D10").Value = 562
$ws.Range("E10").Value = 563
$ws.Range("F10").Value = 564
$ws.Range("C11").Value = 571
$ws.Range("D11").Value = 572
$ws.Range("E11").Value = 573
$ws.Range("F11").Value = 574
# restore/align formatting for Station5
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Style = "Normal"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Style = "Normal"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Style = "Normal"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Style = "Normal"

# ===== restore sheet selections (order matters: last Select() wins the active tab) =====
$ws = $wb.Worksheets.Item("Station1")
$ws.Select()
$ws.Range("I8").Select()

$ws = $wb.Worksheets.Item("Station2")
$ws.Select()
$ws.Range("F2:G11").Select()

$ws = $wb.Worksheets.Item("Station3")
$ws.Select()
$ws.Range("C2:D11").Select()

$ws = $wb.Worksheets.Item("Station4")
$ws.Select()
$ws.Range("C2:E11").Select()

$ws = $wb.Worksheets.Item("Station5")
$ws.Select()
$ws.Range("I6").Select()
